## keck added as an admin
## - New sprint sheet "2017.05.12" added (copied from the last sprint sheet
##   template), last week's sheet "2017.28.11" gets its actual "Worked" hours
##   filled in, and the Summary log gets a new row pointing at the new sprint.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet: append the new sprint log row (do this FIRST so the new
#    shared string "Release" is interned before the strings used on the new
#    sprint sheet, matching shared-string insertion order).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

# Clone the date-format/style from the prior row before writing the values,
# so the new date cell keeps the same number format as the rows above it.
$summary.Range("A20").Copy()
$summary.Range("A21").PasteSpecial(-4122)   # xlPasteFormats
$summary.Range("B20").Copy()
$summary.Range("B21").PasteSpecial(-4122)   # xlPasteFormats

$summary.Range("A21").Value = 43074
$summary.Range("B21").Value = "Release"
$summary.Activate()
$summary.Range("B21").Select()

# ---------------------------------------------------------------------------
# 2. Create the new sprint sheet "2017.05.12" by duplicating the previous
#    week's sheet (keeps layout/formulas/column widths identical), then
#    clear it back to a fresh sprint (new tasks/estimates, nothing worked
#    yet).
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item("2017.28.11")
$last.Copy([System.Reflection.Missing]::Value, $last)
$newSheet = $wb.Worksheets.Item($last.Index + 1)
$newSheet.Name = "2017.05.12"

$newSheet.Range("B3").Value = "Release and Planning"
$newSheet.Range("C3").Value = 6
$newSheet.Range("D3").Value = 6
$newSheet.Range("E3").Value = 0

$newSheet.Range("B8").Value = "Planning"
$newSheet.Range("C8").Value = 4
$newSheet.Range("D8").Value = 4
$newSheet.Range("E8").Value = 0

$newSheet.Range("C9").Value = 2
$newSheet.Range("D9").Value = 2
$newSheet.Range("E9").Value = 0

$newSheet.Columns.Item(2).ColumnWidth = 18.8

$newSheet.Range("E9").Select()

# ---------------------------------------------------------------------------
# 3. Fill in the actual "Worked" hours (column E) on last week's sheet,
#    "2017.28.11" -- this ripples through the existing D-E / SUM formulas.
# ---------------------------------------------------------------------------
$last.Range("E3").Value = 6
$last.Range("E8").Value = 3
$last.Range("E9").Value = 6
$last.Range("A1:F17").Select()

# ---------------------------------------------------------------------------
# 4. Leave the new sprint sheet as the active tab, scrolled so it is visible.
# ---------------------------------------------------------------------------
$newSheet.Activate()
$wb.Windows.Item(1).ScrollWorkbookTabs(10)
